# chore: JP tenth anni
# Updates the "servents" sheet:
#  - swap rows 2 & 3 (sort fix for Foreigner/Berserker)
#  - rename class "Grandcaster" -> "Loregrandcaster" (row 86)
#  - backfill Traditional-Chinese (name_TW) for rows 379-396
#  - backfill Simplified-Chinese (name_CN) for rows 406-420
#  - append 12 new servants as rows 436-447

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 2 & 3 swap -------------------------------------------------
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 'Berserker'
$ws.Range("D2").Value = 'ヘンリー・ジキル＆ハイド'
$ws.Range("E2").Value = '亨利·傑基爾＆海德'
$ws.Range("F2").Value = '亨利．傑基爾＆海德'

$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 'Foreigner'
$ws.Range("D3").Value = '蒼崎青子'
$ws.Range("E3").Value = '蒼崎青子'
$ws.Range("F3").Value = ""

# --- row 86: className typo/rename ------------------------------------
$ws.Range("C86").Value = 'Loregrandcaster'

# --- backfill name_TW (column F) for rows 379-396 ---------------------
$ws.Range("F379").Value = '幼體／迪亞馬特'
$ws.Range("F380").Value = '所多瑪之獸／德拉科'
$ws.Range("F381").Value = '洛庫斯塔'
$ws.Range("F382").Value = '瑟坦特'
$ws.Range("F383").Value = '果心居士'
$ws.Range("F384").Value = '怖軍'
$ws.Range("F385").Value = '難敵'
$ws.Range("F386").Value = '杜爾迦'
$ws.Range("F387").Value = '美杜莎'
$ws.Range("F388").Value = '雨之魔女托內莉可'
$ws.Range("F389").Value = '阿爾托莉亞．Caster'
$ws.Range("F390").Value = '鈴鹿御前〔暑假〕'
$ws.Range("F391").Value = '克洛伊．馮．愛因茲貝倫'
$ws.Range("F392").Value = '諾克娜蕾雅．雅蘭杜'
$ws.Range("F393").Value = '美露莘'
$ws.Range("F394").Value = 'ＵＤＫ－巴格斯特'
$ws.Range("F395").Value = '凱特．庫．米可科爾'
$ws.Range("F396").Value = '旺吉娜'

# --- backfill name_CN (column E) for rows 406-420 ---------------------
$ws.Range("E406").Value = '源賴光／醜御前'
$ws.Range("E407").Value = '由井正雪'
$ws.Range("E408").Value = '宮本伊織'
$ws.Range("E409").Value = '安德洛墨達'
$ws.Range("E410").Value = '瑪麗·安託瓦內特〔Alter〕'
$ws.Range("E411").Value = '耀星哈桑'
$ws.Range("E412").Value = '巖窟王　基督山'
$ws.Range("E413").Value = '亞歷山德羅·迪·卡利奧斯特羅'
$ws.Range("E414").Value = 'Ｅ－火瑪麗'
$ws.Range("E415").Value = 'Ｅ－水瑪麗'
$ws.Range("E416").Value = '蒼崎青子'
$ws.Range("E417").Value = '靜希草十郎'
$ws.Range("E418").Value = '久遠寺有珠'
$ws.Range("E419").Value = '響＆千鍵'
$ws.Range("E420").Value = '埃列什基伽勒'

# --- append new servants: rows 436-447 ---------------------------------
$ws.Range("A436").Value = 433
$ws.Range("B436").Value = 5
$ws.Range("C436").Value = 'Lancer'
$ws.Range("D436").Value = 'ビショーネ'
$ws.Range("E436").Value = ""
$ws.Range("F436").Value = ""

$ws.Range("A437").Value = 434
$ws.Range("B437").Value = 4
$ws.Range("C437").Value = 'Saber'
$ws.Range("D437").Value = '黒姫'
$ws.Range("E437").Value = ""
$ws.Range("F437").Value = ""

$ws.Range("A438").Value = 435
$ws.Range("B438").Value = 5
$ws.Range("C438").Value = 'Caster'
$ws.Range("D438").Value = '小野小町'
$ws.Range("E438").Value = ""
$ws.Range("F438").Value = ""

$ws.Range("A439").Value = 436
$ws.Range("B439").Value = 0
$ws.Range("C439").Value = 'Uolgamariegrandcollection'
$ws.Range("D439").Value = 'Ｅ－グランマリー'
$ws.Range("E439").Value = ""
$ws.Range("F439").Value = ""

$ws.Range("A440").Value = 437
$ws.Range("B440").Value = 5
$ws.Range("C440").Value = 'Pretender'
$ws.Range("D440").Value = 'ダンテ・アリギエーリ'
$ws.Range("E440").Value = ""
$ws.Range("F440").Value = ""

$ws.Range("A441").Value = 438
$ws.Range("B441").Value = 5
$ws.Range("C441").Value = 'Ruler'
$ws.Range("D441").Value = 'メタトロン・ジャンヌ'
$ws.Range("E441").Value = ""
$ws.Range("F441").Value = ""

$ws.Range("A442").Value = 439
$ws.Range("B442").Value = 4
$ws.Range("C442").Value = 'Lancer'
$ws.Range("D442").Value = 'アショカ王'
$ws.Range("E442").Value = ""
$ws.Range("F442").Value = ""

$ws.Range("A443").Value = 440
$ws.Range("B443").Value = 5
$ws.Range("C443").Value = 'Berserker'
$ws.Range("D443").Value = 'リリス'
$ws.Range("E443").Value = ""
$ws.Range("F443").Value = ""

$ws.Range("A444").Value = 441
$ws.Range("B444").Value = 5
$ws.Range("C444").Value = 'Pretender'
$ws.Range("D444").Value = 'テュフォン・エフェメロス'
$ws.Range("E444").Value = ""
$ws.Range("F444").Value = ""

$ws.Range("A445").Value = 442
$ws.Range("B445").Value = 5
$ws.Range("C445").Value = 'Lancer'
$ws.Range("D445").Value = 'インドラ'
$ws.Range("E445").Value = ""
$ws.Range("F445").Value = ""

$ws.Range("A446").Value = 443
$ws.Range("B446").Value = 0
$ws.Range("C446").Value = 'Uolgamariestellarcollection'
$ws.Range("D446").Value = 'Ｅ－ステラマリー'
$ws.Range("E446").Value = ""
$ws.Range("F446").Value = ""

# row 447's className is the literal text "40" (not a number) in the
# source data, so force the cell to Text before writing it.
$ws.Range("A447").Value = 444
$ws.Range("B447").Value = 5
$ws.Range("C447").NumberFormat = "@"
$ws.Range("C447").Value = '40'
$ws.Range("D447").Value = 'Ｕ－オルガマリー'
$ws.Range("E447").Value = ""
$ws.Range("F447").Value = ""
